# Auto-generated edit script: updates cryptos price/volume columns
# per the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.080.18"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "2.573.78"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("E4").Value = "  +0.04%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "572.65"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +2.53%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "142.79"
$cell.ClearFormats()
$ws.Range("E6").Value = "  -0.76%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.ClearFormats()
$ws.Range("E7").Value = "  +0.08%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.595"
$cell.ClearFormats()
$ws.Range("E8").Value = "  -0.42%  "
$ws.Range("D9").Value = "2.581.24"
$ws.Range("E9").Value = "  -1.59%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "6.69"
$cell.ClearFormats()
$ws.Range("E10").Value = "  -1.82%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.103"
$cell.ClearFormats()
$ws.Range("E11").Value = "  +2.56%  "
$ws.Range("E12").Value = "  +11.63%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.344"
$cell.ClearFormats()
$ws.Range("E13").Value = "  +2.32%  "
$ws.Range("D14").Value = "3.030.87"
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("D15").Value = "59.136.39"
$ws.Range("E15").Value = "  +0.24%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "22.31"
$cell.ClearFormats()
$ws.Range("E16").Value = "  +5.59%  "
$ws.Range("E17").Value = "  +3.11%  "
$ws.Range("D18").Value = "2.585.47"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("E19").Value = "  +1.50%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "335.76"
$cell.ClearFormats()
$ws.Range("E20").Value = "  -0.69%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "10.26"
$cell.ClearFormats()
$ws.Range("E21").Value = "  +1.30%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "6.27"
$cell.ClearFormats()
$ws.Range("E22").Value = "  +1.43%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.ClearFormats()
$ws.Range("E23").Value = "  -0.09%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "64.55"
$cell.ClearFormats()
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "0.464"
$cell.ClearFormats()
$ws.Range("E25").Value = "  +8.17%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "1.01"
$cell.ClearFormats()
$ws.Range("E26").Value = "  +0.96%  "
$ws.Range("E27").Value = "  -1.29%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "7.27"
$cell.ClearFormats()
$ws.Range("E28").Value = "  +1.11%  "
$ws.Range("D29").Value = "0.0₃0780"
$ws.Range("E29").Value = "  +2.12%  "
$ws.Range("E30").Value = "  +0.06%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.68"
$cell.ClearFormats()
$ws.Range("E31").Value = "  +0.06%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "159.62"
$cell.ClearFormats()
$ws.Range("E32").Value = "  +3.13%  "
$ws.Range("E33").Value = "  +0.41%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "18.95"
$cell.ClearFormats()
$ws.Range("E34").Value = "  +0.16%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "4.02"
$cell.ClearFormats()
$ws.Range("E35").Value = "  +0.91%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "1.16"
$cell.ClearFormats()
$ws.Range("E36").Value = "  +2.40%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.872"
$cell.ClearFormats()
$ws.Range("E37").Value = "  -4.63%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.876"
$cell.ClearFormats()
$ws.Range("E38").Value = "  -4.49%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "37.36"
$cell.ClearFormats()
$ws.Range("E39").Value = "  +0.51%  "
$ws.Range("E40").Value = "  +0.87%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "295.74"
$cell.ClearFormats()
$ws.Range("E41").Value = "  +3.43%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "3.66"
$cell.ClearFormats()
$ws.Range("E42").Value = "  +1.46%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.ClearFormats()
$ws.Range("E43").Value = "  +0.33%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "130.72"
$cell.ClearFormats()
$ws.Range("E44").Value = "  +10.47%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.0976"
$cell.ClearFormats()
$ws.Range("E45").Value = "  +2.18%  "
$ws.Range("E46").Value = "  -1.04%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.0536"
$cell.ClearFormats()
$ws.Range("E47").Value = "  -0.66%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "10.65"
$cell.ClearFormats()
$ws.Range("E48").Value = "  +0.35%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "19.16"
$cell.ClearFormats()
$ws.Range("E49").Value = "  +1.80%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.0233"
$cell.ClearFormats()
$ws.Range("E50").Value = "  +1.90%  "
$ws.Range("D51").Value = "1.950.50"
$ws.Range("E51").Value = "  +0.06%  "
